# Going through the dataset, updating.
#
# This pass:
#   1) Adds a new "Other found locations" column (I) that records which
#      extra source(s) (PMC / PMC+Springer / none) were found for each
#      reference while re-checking the dataset.
#   2) Normalises a handful of "Authors" (column E) values whose delimiter
#      spacing / encoding were cleaned up on this re-check.
#   3) Flips a few rows whose "ID" (F) / "ID Format" (G) could not actually
#      be resolved back to "not found" / "N/A".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) New column header -------------------------------------------------
$ws.Range("I1").Value = "Other found locations"

# --- 2) Refreshed "Authors" (E) values ------------------------------------
$ws.Range("E2").Value  = "[Sapna R.%Kudchadkar%NULL%1,   Christopher L.%Carroll%NULL%1]"
$ws.Range("E3").Value  = "[Sohaib R%Rufai%sohaibrufai@gmail.com%1,   Catey%Bunce%Catey.bunce@kcl.ac.uk%2,   Catey%Bunce%Catey.bunce@kcl.ac.uk%0]"
$ws.Range("E4").Value  = "[Travis%Sanchez%NULL%2,   Kamal%Al Nasr%NULL%2,   Kamal%Al Nasr%NULL%0,   Ross%Gore%NULL%1,   Abdullah%Wahbeh%NULL%2,   Abdullah%Wahbeh%NULL%0,   Tareq%Nasralah%t.nasralah@northeastern.edu%2,   Tareq%Nasralah%t.nasralah@northeastern.edu%0,   Mohammad%Al-Ramahi%NULL%2,   Mohammad%Al-Ramahi%NULL%0,   Omar%El-Gayar%NULL%2,   Omar%El-Gayar%NULL%0]"
$ws.Range("E9").Value  = "[Vittorio%Gebbia%NULL%1,   Dario%Piazza%NULL%1,   Maria Rosaria%Valerio%NULL%1,   Nicolò%Borsellino%NULL%1,   Alberto%Firenze%NULL%1]"
$ws.Range("E11").Value = "[Maria Renee%Jimenez‐Sotomayor%NULL%1,   Carolina%Gomez‐Moreno%NULL%1,   Enrique%Soto‐Perez‐de‐Celis%enrique.sotop@incmnsz.mx%1]"
$ws.Range("E12").Value = "[Greg%Kawchuk%greg.kawchuk@ualberta.ca%1,   Jan%Hartvigsen%jhartvigsen@health.sdu.dk%2,   Jan%Hartvigsen%jhartvigsen@health.sdu.dk%0,   Stan%Innes%S.Innes@murdoch.edu.au%1,   J. Keith%Simpson%k.simpson@murdoch.edu.au%1,   Brian%Gushaty%bgushaty@gushaty.com%1]"
$ws.Range("E13").Value = "[Gunther%Eysenbach%NULL%0,   Jon-Patrick%Allem%NULL%2,   Jon-Patrick%Allem%NULL%0,   Richard%Zowalla%NULL%1,   Wasim%Ahmed%Wasim.Ahmed@Newcastle.ac.uk%2,   Wasim%Ahmed%Wasim.Ahmed@Newcastle.ac.uk%0,   Josep%Vidal-Alaball%NULL%2,   Josep%Vidal-Alaball%NULL%0,   Joseph%Downing%NULL%2,   Joseph%Downing%NULL%0,   Francesc%López Seguí%NULL%2,   Francesc%López Seguí%NULL%0]"
$ws.Range("E16").Value = "[Servet%Aker%servetaker@gmail.com%1,   Özlem%Mıdık%NULL%2,   Özlem%Mıdık%NULL%0]"
$ws.Range("E18").Value = "[Flecha%Ramón%coreGivesNoEmail%1,  Guo%Mengna%coreGivesNoEmail%1,  Pulido%Rodríguez Cristina%coreGivesNoEmail%1,  Ramis%Salas Maria del Mar%coreGivesNoEmail%1,  Redondo%Sama Gisela%coreGivesNoEmail%1,  Villarejo%Carballido Beatriz%coreGivesNoEmail%1]"
$ws.Range("E19").Value = "[Bahar%Yuksel%baharyl86@gmail.com%1,   Kubra%Cakmak%NULL%1]"
$ws.Range("E20").Value = "[Gunther%Eysenbach%NULL%0,   Nazakat%Hamassed%NULL%2,   Nazakat%Hamassed%NULL%0,   Hardawan%Kakashekh%NULL%1,   Muhammad%Saud%NULL%1,   Mohammad Amin%Bahrami%NULL%1,   Araz Ramazan%Ahmad%araz.ahmad85@uor.edu.krd%2,   Araz Ramazan%Ahmad%araz.ahmad85@uor.edu.krd%0,   Hersh Rasool%Murad%NULL%2,   Hersh Rasool%Murad%NULL%0]"
$ws.Range("E21").Value = "[Alexander%Muacevic%NULL%0,   John R%Adler%NULL%0,   Ramez%Kouzy%NULL%2,   Ramez%Kouzy%NULL%0,   Joseph%Abi Jaoude%NULL%1,   Afif%Kraitem%NULL%1,   Molly B%El Alam%NULL%1,   Basil%Karam%NULL%1,   Elio%Adib%NULL%1,   Jabra%Zarka%NULL%1,   Cindy%Traboulsi%NULL%1,   Elie W%Akl%NULL%1,   Khalil%Baddour%NULL%1]"

# --- 3) Rows whose ID / ID Format reverted to "not found" / "N/A" --------
$ws.Range("F7").Value  = "not found"
$ws.Range("G7").Value  = "N/A"
$ws.Range("F14").Value = "not found"
$ws.Range("G14").Value = "N/A"
$ws.Range("F18").Value = "not found"
$ws.Range("G18").Value = "N/A"

# --- 1b) Per-row "Other found locations" values ---------------------------
# Rows whose reference was additionally located on PubMed Central:
$ws.Range("I2").Value  = "_PMC"
$ws.Range("I3").Value  = "_PMC"
$ws.Range("I4").Value  = "_PMC"
$ws.Range("I5").Value  = "_PMC"
$ws.Range("I9").Value  = "_PMC"
$ws.Range("I11").Value = "_PMC"
$ws.Range("I13").Value = "_PMC"
$ws.Range("I19").Value = "_PMC"
$ws.Range("I20").Value = "_PMC"
$ws.Range("I21").Value = "_PMC"

# Rows whose reference was located on PubMed Central as well as Springer:
$ws.Range("I12").Value = "_PMC_Springer"
$ws.Range("I16").Value = "_PMC_Springer"

# Rows 6, 7, 8, 10, 14, 15, 17 and 18 record "no other location found",
# i.e. an empty string in the "Other found locations" column -- which,
# same as in the Excel UI, is simply a blank cell.
$blankRows = @(6, 7, 8, 10, 14, 15, 17, 18)
foreach ($r in $blankRows) {
    $ws.Range("I$r").Value = ""
}
